$d = $word.ActiveDocument

# The styles below have a <w:rPr> child-element order that does not match
# the WordprocessingML schema (color must come *after* b/i, not before).
# Re-asserting the Bold/Italic flag on each style's Font forces the engine
# to rewrite rPr in schema-compliant order (b, i, ..., color), fixing the
# OOXMLValidator "Sch_UnexpectedElementContentExpectingComplex" warning.

$boldStyles = @(
    "KeywordTok",
    "ImportTok",
    "AnnotationTok",
    "CommentVarTok",
    "ControlFlowTok",
    "InformationTok",
    "WarningTok",
    "AlertTok",
    "ErrorTok"
)

$italicStyles = @(
    "CommentTok",
    "DocumentationTok",
    "AnnotationTok",
    "CommentVarTok",
    "InformationTok",
    "WarningTok"
)

foreach ($styleName in $boldStyles) {
    $style = $d.Styles($styleName)
    $style.Font.Bold = $true
}

foreach ($styleName in $italicStyles) {
    $style = $d.Styles($styleName)
    $style.Font.Italic = $true
}
